# Generate Report for Handback
# Update timestamp strings on the three worksheets to reflect the latest
# handoff/handback run times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file (row 2)
$wsOverview.Range("G2").Value = "2016-08-17 23:02:58"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2)
$wsZhCn.Range("H2").Value = "2016-08-17 23:02:53"
$wsZhCn.Range("K2").Value = "2016-08-17 23:03:17"

# de-de sheet: Correspond Handoff Datetime (row 2) shares the same value as
# Overview!G2, and Correspond Handback DateTime (row 2) gets a new value.
$wsDeDe.Range("H2").Value = "2016-08-17 23:02:58"
$wsDeDe.Range("K2").Value = "2016-08-17 23:03:24"
